$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21 (pushes existing rows 21..133 down to 22..134)
$ws.Rows("21:21").Insert()

# Populate the newly inserted row 21 with the new weekly record.
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(21, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(21, 4).Value = 45107
$ws.Cells.Item(21, 5).Value = 15
$ws.Cells.Item(21, 6).Value = 100112038
$ws.Cells.Item(21, 7).Value = 'Cebollín baby'
$ws.Cells.Item(21, 8).Value = 'Sin especificar'
$ws.Cells.Item(21, 9).Value = 'Primera'
$ws.Cells.Item(21, 10).Value = 300
$ws.Cells.Item(21, 11).Value = 1800
$ws.Cells.Item(21, 12).Value = 2000
$ws.Cells.Item(21, 13).Value = 1900
$ws.Cells.Item(21, 14).Value = '$/paquete 1,5 a 2 kilos'
$ws.Cells.Item(21, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(21, 16).Value = 950
$ws.Cells.Item(21, 17).Value = 2
$ws.Cells.Item(21, 18).Value = 'Hortaliza'

# Match the date formatting used by the rest of column D.
$ws.Cells.Item(21, 4).NumberFormat = $ws.Cells.Item(22, 4).NumberFormat
